# Adapt column header formatting to respective input file names:
#   "<name>_old" -> "<name>_FV2310"
#   "<name>_new" -> "<name>_FV2404"
# then wrap the data range in a table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1) suffixes -----------------------------
$lastCol = $ws.UsedRange.Columns.Count
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $header = $cell.Value2
    if ($header -ne $null) {
        $newHeader = $header
        if ($header.EndsWith("_old")) {
            $newHeader = $header.Substring(0, $header.Length - 4) + "_FV2310"
        } elseif ($header.EndsWith("_new")) {
            $newHeader = $header.Substring(0, $header.Length - 4) + "_FV2404"
        }
        if ($newHeader -ne $header) {
            $cell.Value2 = $newHeader
        }
    }
}

# --- 2. Turn the used range into a real Excel table -------------------------
$lastRow = $ws.UsedRange.Rows.Count
$dataRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row -----------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
